# B6-PowerPoint.pptx edit:
#   1. Three tables (one each on slides 14, 15, 16) get their table style
#      swapped from {D9246204-775E-4A70-9B95-437E35A5998C} to
#      {57A20238-2F8F-4438-8F75-A78454F987ED}.
#   2. The deck's theme palette ("Integral" / Red Violet) is swapped back to
#      the stock "Office Theme" palette on the slide master's theme part.

$p = $ppt.ActivePresentation

# --- 1. Retarget the table styles -----------------------------------------
$newTableStyle = "{57A20238-2F8F-4438-8F75-A78454F987ED}"

for ($idx = 1; $idx -le $p.Slides.Count; $idx++) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Swap the theme colour scheme back to the Office Theme palette -----
$slide1 = $p.Slides.Item(1)
$colorScheme = $slide1.ThemeColorScheme

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as COM RGB() longs (0x00BBGGRR) to match VBA's ColorFormat.RGB.
$officeTheme = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = $officeTheme[$i - 1]
}
